# Actualización automática: refresh the "urbana" slug columns in row 2
# so the measure-slug values match the dash style used elsewhere in the
# sheet (single dash -> triple dash), e.g. "urbana-cuota-integra-euros"
# becomes "urbana---cuota-integra-euros".
#
# Affected cells (row 2 holds the slugified header/measure-id row):
#   B2 -> urbana---cuota-integra-euros
#   F2 -> urbana---numero-de-recibos
#   G2 -> urbana---base-imponible-miles-de-euros
#   K2 -> urbana---cuota-liquida-euros
#   N2 -> urbana---base-liquidable-miles-de-euros

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "urbana---cuota-integra-euros"
$ws.Range("F2").Value = "urbana---numero-de-recibos"
$ws.Range("G2").Value = "urbana---base-imponible-miles-de-euros"
$ws.Range("K2").Value = "urbana---cuota-liquida-euros"
$ws.Range("N2").Value = "urbana---base-liquidable-miles-de-euros"

# The trailing O1:P1 / O2:O5 cells carry no content (only leftover
# styling from the original export) and are dropped from the sheet on
# this refresh.
$ws.Range("O1:P1").Clear()
$ws.Range("O2:O5").Clear()
